$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Title paragraph: "Liberal Arts " + [bookmark] + "Education" + " Reflection"
#    becomes "Liberal Arts Education" + " Reflection" (bookmark removed here,
#    it moves down into the body paragraph below).
# ------------------------------------------------------------------
$d.Bookmarks.Item("_GoBack").Delete()
$d.Content.Find.Execute("Liberal Arts Education", $true, $false, $false, $false, $false, $true, 1, $false, "Liberal Arts Education", 2) | Out-Null

# ------------------------------------------------------------------
# 2. Body paragraph formatting: add <w:ind w:right="-90"/>
#    (-90 twips = -4.5 points)
# ------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$p2.Range.ParagraphFormat.RightIndent = -4.5

# ------------------------------------------------------------------
# 3. Fix typo: "a liberal educations means" -> "a liberal education means"
# ------------------------------------------------------------------
$d.Content.Find.Execute("a liberal educations means", $true, $false, $false, $false, $false, $true, 1, $false, "a liberal education means", 2) | Out-Null

# ------------------------------------------------------------------
# 4. Append the new closing sentence to the body paragraph, right
#    before its paragraph mark, and give it the same sz=24 run
#    formatting as the rest of the paragraph.
# ------------------------------------------------------------------
$insertPos = $p2.Range.End - 1
$newText = " Please provide some depth to your response (i.e., use all of the space provided)."
$insertRng = $d.Range($insertPos, $insertPos)
$insertRng.InsertAfter($newText)
$newRng = $d.Range($insertPos, $insertPos + $newText.Length)
$newRng.Font.Size = 12

# ------------------------------------------------------------------
# 5. Re-insert the _GoBack bookmark at its new home: right after
#    "...your" and right before " response (i.e., ...)".
# ------------------------------------------------------------------
$bmFind = $d.Content
$bmFind.Find.Execute("some depth to your", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmPos = $bmFind.End
$bmRng = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRng)
